$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: E14 gets the new "^" marker, styled like the adjacent Mg/unit cells
# (text number format, matching D14/H14) instead of the previously blank cell.
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "^"

# Row 17: Phosphatidylserine's amount moves from a text "50 mg" in C17 into a
# proper numeric amount (matching C15/C16), and its unit "Mg" moves into D17
# (matching D15/D16), cleaning up the row to be consistent with rows 15-16.
$ws.Range("C17").NumberFormat = "General"
$ws.Range("C17").Value = 50

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "Mg"
